$d = $word.ActiveDocument

# 1) Merge the three runs ("...un " + "refactor" + " en el proyecto...cambiarlos. ")
#    into a single run of plain text, removing the spell-check markers around
#    "refactor". A Find/Replace that spans the run boundaries merges the text
#    and collapses it into a single run (and drops the now-orphaned
#    w:proofErr markers).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "produjo incluso un refactor en el proyecto", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "produjo incluso un refactor en el proyecto", 2) | Out-Null

# 2) Append a brand-new sentence about the final grade to the end of the
#    last paragraph, as its own separate run (not merged into the existing
#    run). Toggling a character property on the freshly-inserted range
#    (and then reverting it) is enough to make the engine keep it as a
#    distinct run while leaving the visible formatting unchanged.
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertAfter("Por ende, aunque mi empeño en el trabajo fue bastante grande, creo que pudo ser por mucho mejorable, y según yo, me quedé corto con algunas funcionalidades. Es por ello, que mi nota en el proyecto es 4.5.")
$endRange.Bold = 1
$endRange.Bold = 0
